$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (was date 44340, now date 44372)
$ws.Range("D2").Value = 44372
$ws.Range("O2").Value = 'Región Metropolitana'

# Row 3 (was date 44369, now date 44371)
$ws.Range("D3").Value = 44371
$ws.Range("J3").Value = 150
$ws.Range("K3").Value = 6500
$ws.Range("L3").Value = 6500
$ws.Range("M3").Value = 6500
$ws.Range("N3").Value = '$/caja 36 atados'
$ws.Range("P3").Value = 181
$ws.Range("Q3").Value = 36

# Row 4 (was date 44342, now date 44364)
$ws.Range("D4").Value = 44364
$ws.Range("J4").Value = 100
$ws.Range("O4").Value = 'Región Metropolitana'

# Row 5 (was date 44357, now date 44340)
$ws.Range("D5").Value = 44340
$ws.Range("K5").Value = 7000
$ws.Range("L5").Value = 7000
$ws.Range("M5").Value = 7000
$ws.Range("N5").Value = '$/caja 36 atados'
$ws.Range("P5").Value = 194
$ws.Range("Q5").Value = 36

# Row 6 (was date 44355, now date 44342)
$ws.Range("D6").Value = 44342
$ws.Range("O6").Value = 'Región del Maule'

# Row 7 (was date 44348, now date 44376)
$ws.Range("D7").Value = 44376
$ws.Range("K7").Value = 6500
$ws.Range("L7").Value = 6500
$ws.Range("M7").Value = 6500
$ws.Range("O7").Value = 'Región Metropolitana'
$ws.Range("P7").Value = 181

# Row 8 (was date 44358, now date 44369)
$ws.Range("D8").Value = 44369
$ws.Range("J8").Value = 100
$ws.Range("N8").Value = '$/caja 20 docenas'
$ws.Range("P8").Value = 7000
$ws.Range("Q8").Value = 1

# Row 9 (was date 44362, now date 44354)
$ws.Range("D9").Value = 44354
$ws.Range("J9").Value = 150
$ws.Range("K9").Value = 7000
$ws.Range("L9").Value = 7000
$ws.Range("M9").Value = 7000
$ws.Range("O9").Value = 'Región del Maule'
$ws.Range("P9").Value = 194

# Row 10 (was date 44364, now date 44386)
$ws.Range("D10").Value = 44386
$ws.Range("J10").Value = 200
$ws.Range("K10").Value = 6500
$ws.Range("L10").Value = 6500
$ws.Range("M10").Value = 6500
$ws.Range("P10").Value = 181

# Row 11 (was date 44376, now date 44358)
$ws.Range("D11").Value = 44358
$ws.Range("K11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("M11").Value = 7000
$ws.Range("P11").Value = 194

# Row 12 (was date 44386, now date 44348)
$ws.Range("D12").Value = 44348
$ws.Range("J12").Value = 150
$ws.Range("K12").Value = 7000
$ws.Range("L12").Value = 7000
$ws.Range("M12").Value = 7000
$ws.Range("O12").Value = 'Región del Maule'
$ws.Range("P12").Value = 194

# Row 13 (was date 44371, now date 44362)
$ws.Range("D13").Value = 44362
$ws.Range("J13").Value = 100

# Row 14 (was date 44354, now date 44357)
$ws.Range("D14").Value = 44357
$ws.Range("K14").Value = 6500
$ws.Range("L14").Value = 6500
$ws.Range("M14").Value = 6500
$ws.Range("N14").Value = '$/caja 20 docenas'
$ws.Range("P14").Value = 6500
$ws.Range("Q14").Value = 1

# Row 15 (was date 44372, now date 44355)
$ws.Range("D15").Value = 44355
